$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in marks for students who previously had blank scores ---

# Row 4
$ws.Range("C4").Value = 9
$ws.Range("D4").Value = 8
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 10

# Row 8 (total becomes an explicit SUM formula)
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 10
$ws.Range("G8").Formula = "=SUM(B8:F8)"

# Row 12
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 7
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 10

# Row 14
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 4
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 10

# --- Row heights settle to 12.75pt (rows 2-17) after the edits ---
$ws.Rows("2:17").RowHeight = 12.75

# --- Move the active selection to A15 ---
$ws.Range("A15").Select()
